$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing header row down to row 5, freeing rows 1-4 for the new
# report-parameter block (Sales Agent Name / From Date / To Date).
$ws.Rows("1:4").Insert()

# Row 1: "Sales Agent Name:" label (bold, no number format -> new style index 4)
$ws.Range("A1").Value = "Sales Agent Name:"
$ws.Range("A1").Font.Bold = $true

# Row 2: "From Date:" label + adjoining date-entry cell (date format -> new style index 5)
$ws.Range("A2").Value = "From Date:"
$ws.Range("A2").Font.Bold = $true
$ws.Range("B2").NumberFormat = "mm-dd-yy"

# Row 3: "To Date:" label + adjoining date-entry cell (reuse styles 4 and 5)
$ws.Range("A3").Value = "To Date:"
$ws.Range("A3").Font.Bold = $true
$ws.Range("A1").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)

# Columns D and F carry an explicit column-level style; rows 1-3 get explicit
# "Normal"-styled placeholder cells there so they don't inherit it.
$ws.Range("D1").Style = "Normal"
$ws.Range("F1").Style = "Normal"
$ws.Range("D2").Style = "Normal"
$ws.Range("F2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("F3").Style = "Normal"

# Create the non-bold "0.00" number style (new style index 6) via a throwaway
# cell below the table, then drop the cell again - only the style survives,
# ready to be the default formatting for future data rows under column H.
$ws.Range("H6").NumberFormat = "0.00"
$ws.Range("H6").Delete()

# Row 5 (former row 1): two extra header cells for the new columns.
$ws.Range("G5").Value = "Sold Stocks"
$ws.Range("A1").Copy()
$ws.Range("G5").PasteSpecial(-4122)

$ws.Range("H5").Value = "Sold Stocks Value"
$ws.Range("H5").Font.Bold = $true
$ws.Range("H5").NumberFormat = "0.00"

[void]$ws.Range("H6").Select()
